# Update parametric survival model parameter estimates/SEs and their
# covariance matrices for each distribution fitted to the Falkson 1991
# OS DTIC arm, as part of creating the stacked BUGS data set used for
# the multivariate NMA.

$wb = $excel.ActiveWorkbook

# --- weibull (shape / scale) ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.76488446758512
$ws.Range("C2").Value = 0.2959383108534
$ws.Range("B3").Value = 0.174964827340232
$ws.Range("C3").Value = 0.114840037842542

# --- lognormal (meanlog / sdlog) ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.98537693246169
$ws.Range("C2").Value = 0.309598793339315
$ws.Range("B3").Value = -0.971956116966394
$ws.Range("C3").Value = 0.12666401900976

# --- llogis (shape / scale) ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.07102831617382
$ws.Range("C2").Value = 0.21136996067186
$ws.Range("B3").Value = 1.68163101716978
$ws.Range("C3").Value = 0.224759741543284

# --- gompertz (shape / rate) ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.54933555186575
$ws.Range("C2").Value = 0.249008871886208
$ws.Range("B3").Value = 0.0102466828139336
$ws.Range("C3").Value = 0.0195441948683702

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0875794838307637
$ws.Range("B2").Value = -0.0198566769345156
$ws.Range("A3").Value = -0.0198566769345156
$ws.Range("B3").Value = 0.0131882342916765

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0958514128371598
$ws.Range("B2").Value = -0.0289126257453649
$ws.Range("A3").Value = -0.0289126257453649
$ws.Range("B3").Value = 0.0160437737117049

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0446772602744237
$ws.Range("B2").Value = 0.00891023701979488
$ws.Range("A3").Value = 0.00891023701979488
$ws.Range("B3").Value = 0.0505169414186039

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0620054182780418
$ws.Range("B2").Value = -0.00214364487714964
$ws.Range("A3").Value = -0.00214364487714964
$ws.Range("B3").Value = 0.000381975553052829
